# Apply cryptos.xlsx price/volume update (GitHub Actions scheduled refresh)
# Commit: Updated cryptos list on Sat Feb 24 18:34:39 UTC 2024 with GitHub Actions
#
# Notes:
#  - Column D holds price text that sometimes looks like a plain number
#    (e.g. "0.593", "382.87"). Excel auto-converts such strings to real
#    numbers on assignment, which would corrupt values like "51.569.69"
#    style (thousand-dot) text or introduce floating point noise.
#    For the handful of cells where the new text is a plain decimal number,
#    we prefix the value with a leading apostrophe - exactly what a user
#    does in Excel to force text storage - so the cell keeps storing the
#    exact original string instead of being re-typed as a Number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.569.69"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "2.994.70"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'382.87"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").Value = "'103.58"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.593"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "3.469.14"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'18.43"
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'7.81"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").Value = "3.009.91"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").Value = "'11.14"
$ws.Range("E17").Value = "  +4.21%  "
$ws.Range("D19").Value = "51.586.15"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("D21").Value = "'12.63"
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'70.56"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "'267.83"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'3.22"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("D26").Value = "'7.85"
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").Value = "'10.36"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("D33").Value = "'34.79"
$ws.Range("E33").Value = "  +4.42%  "
$ws.Range("D34").Value = "'51.58"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").Value = "'2.07"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").Value = "'0.0442"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("D39").Value = "'16.81"
$ws.Range("E39").Value = "  +3.54%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.57"
$ws.Range("E40").Value = "  +3.46%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("D43").Value = "'125.01"
$ws.Range("E43").Value = "  +3.98%  "
$ws.Range("E44").Value = "  +9.76%  "
$ws.Range("D45").Value = "'21.51"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").Value = "'2.04"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("D48").Value = "'0.272"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").Value = "2.047.47"
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("D51").Value = "'0.541"
$ws.Range("E51").Value = "  +17.23%  "
